# Updates the cryptos list (prices in column D, 1h volume % in column E)
# to the latest scrape values, and swaps the Kaspa/Aave rows (rank
# positions 38/39) so Aave now ranks just above Kaspa with refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values (e.g. "1.00", "624.34") are digit-only strings that
# Excel would otherwise auto-convert to numbers (dropping trailing zeros /
# using floating point) when assigned via .Value. Force the cell to Text
# format first, assign the literal string, then restore the original
# style so no other formatting changes leak in.
function Set-TextValue($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "76.337.73"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "3.039.23"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws "D5" "200.31"
$ws.Range("E5").Value = "  -1.42%  "
Set-TextValue $ws "D6" "624.34"
$ws.Range("E6").Value = "  +4.67%  "
Set-TextValue $ws "D7" "0.999"
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue $ws "D8" "0.548"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("D10").Value = "3.037.57"
$ws.Range("E10").Value = "  +3.57%  "
Set-TextValue $ws "D11" "0.438"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E13").Value = "  +5.87%  "
$ws.Range("D14").Value = "3.598.35"
$ws.Range("E14").Value = "  +3.55%  "
Set-TextValue $ws "D15" "29.07"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("D16").Value = "76.332.71"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "3.036.10"
$ws.Range("E18").Value = "  +3.34%  "
Set-TextValue $ws "D19" "13.52"
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("E20").Value = "  +1.26%  "
Set-TextValue $ws "D21" "374.91"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  +1.43%  "
Set-TextValue $ws "D24" "73.17"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +1.34%  "
Set-TextValue $ws "D28" "9.80"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  +0.68%  "
Set-TextValue $ws "D30" "1.00"
$ws.Range("E30").Value = "  -0.04%  "
Set-TextValue $ws "D31" "8.27"
$ws.Range("E31").Value = "  +6.01%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("E33").Value = "  +6.15%  "
Set-TextValue $ws "D34" "492.39"
$ws.Range("E34").Value = "  -2.06%  "
Set-TextValue $ws "D35" "0.999"
$ws.Range("E35").Value = "  -0.01%  "
Set-TextValue $ws "D36" "20.64"
$ws.Range("E36").Value = "  +1.74%  "
Set-TextValue $ws "D37" "162.61"
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("E38").Value = "  +2.05%  "
Set-TextValue $ws "D39" "0.383"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D40" "190.33"
$ws.Range("E40").Value = "  +4.71%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D41" "0.115"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  -5.67%  "
$ws.Range("E43").Value = "  +0.01%  "
Set-TextValue $ws "D44" "0.803"
$ws.Range("E44").Value = "  +22.09%  "
Set-TextValue $ws "D45" "5.12"
$ws.Range("E45").Value = "  +2.40%  "
Set-TextValue $ws "D46" "1.26"
$ws.Range("E46").Value = "  +5.20%  "
$ws.Range("E47").Value = "  +4.75%  "
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("E49").Value = "  +4.52%  "
Set-TextValue $ws "D50" "0.604"
$ws.Range("E50").Value = "  +4.24%  "
$ws.Range("E51").Value = "  +4.09%  "
